# Trade #66 closed at 2026-02-17 21:11:44 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#  - Summary sheet: capital / P&L / trade-count roll-up numbers
#  - Strategy Status sheet: MarketMaking strategy row roll-up numbers
#  - All Trades sheet: close out the open MarketMaking trade (row 95) and
#    append a brand-new open MarketMaking trade (row 128)
#  - MarketMaking sheet: same two changes, mirrored into this sheet's own
#    column layout (row 62 closes, row 95 is appended)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.05   # Current Capital
$summary.Range("B4").Value = 0.85      # Total P&L $
$summary.Range("B6").Value = 94        # Total Trades
$summary.Range("B8").Value = 38        # Losing Trades
$summary.Range("B9").Value = 46.81     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.05     # Capital
$status.Range("D5").Value = 61         # Trades
$status.Range("E5").Value = 0.74       # P&L $
$status.Range("F5").Value = 1.05       # P&L %
$status.Range("G5").Value = 49.18      # Win Rate %

# ---------------------------------------------------------------------
# All Trades - columns: A Trade#, B Date, C Time, D Strategy, E Side,
# F Entry Price, G Exit Price, H Status, I P&L%, J P&L$, K Capital After,
# L Exit Reason, M Duration (min), N Entry Slippage, O Exit Slippage,
# P Confidence, Q Entry Reason
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Row 95 (Trade #94) closes out
$allTrades.Cells.Item(95, 7).Value = 0.87        # G95 Exit Price
$allTrades.Cells.Item(95, 8).Value = "CLOSED"    # H95 Status
$allTrades.Cells.Item(95, 9).Value = -1.1364     # I95 P&L %
$allTrades.Cells.Item(95, 10).Value = -0.01      # J95 P&L $
$allTrades.Cells.Item(95, 11).Value = 101.05     # K95 Capital After
$allTrades.Cells.Item(95, 12).Value = "early_exit" # L95 Exit Reason
$allTrades.Cells.Item(95, 13).Value = 0.14       # M95 Duration (min)

# Row 128 (Trade #127) newly appended, still open
$allTrades.Cells.Item(128, 1).Value = 127
$allTrades.Cells.Item(128, 2).NumberFormat = "@"
$allTrades.Cells.Item(128, 2).Value = "2026-02-17"
$allTrades.Cells.Item(128, 3).NumberFormat = "@"
$allTrades.Cells.Item(128, 3).Value = "21:11:37"
$allTrades.Cells.Item(128, 4).Value = "MarketMaking"
$allTrades.Cells.Item(128, 5).Value = "UP"
$allTrades.Cells.Item(128, 6).Value = 0.88
$allTrades.Cells.Item(128, 8).Value = "OPEN"
$allTrades.Cells.Item(128, 9).Value = 0
$allTrades.Cells.Item(128, 10).Value = 0
$allTrades.Cells.Item(128, 11).Value = 101.0646450978375
$allTrades.Cells.Item(128, 13).Value = 0
$allTrades.Cells.Item(128, 14).Value = 0
$allTrades.Cells.Item(128, 15).Value = 0
$allTrades.Cells.Item(128, 16).Value = 0.6
$allTrades.Cells.Item(128, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking - columns: A Trade#, B Date, C Time, D Strategy, E Side,
# F Entry Price, G Exit Price, H Status, I P&L%, J P&L$, K Capital After,
# L Entry Slippage, M Exit Slippage, N Confidence, O Entry Reason,
# P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# Row 62 (Trade #94) closes out
$marketMaking.Cells.Item(62, 7).Value = 0.87        # G62 Exit Price
$marketMaking.Cells.Item(62, 8).Value = "CLOSED"    # H62 Status
$marketMaking.Cells.Item(62, 9).Value = -1.1364     # I62 P&L %
$marketMaking.Cells.Item(62, 10).Value = -0.01      # J62 P&L $
$marketMaking.Cells.Item(62, 11).Value = 101.05     # K62 Capital After
$marketMaking.Cells.Item(62, 16).Value = "early_exit" # P62 Exit Reason
$marketMaking.Cells.Item(62, 17).Value = 0.14       # Q62 Duration (min)

# Row 95 (Trade #127) newly appended, still open
$marketMaking.Cells.Item(95, 1).Value = 127
$marketMaking.Cells.Item(95, 2).NumberFormat = "@"
$marketMaking.Cells.Item(95, 2).Value = "2026-02-17"
$marketMaking.Cells.Item(95, 3).NumberFormat = "@"
$marketMaking.Cells.Item(95, 3).Value = "21:11:37"
$marketMaking.Cells.Item(95, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(95, 5).Value = "UP"
$marketMaking.Cells.Item(95, 6).Value = 0.88
$marketMaking.Cells.Item(95, 8).Value = "OPEN"
$marketMaking.Cells.Item(95, 9).Value = 0
$marketMaking.Cells.Item(95, 10).Value = 0
$marketMaking.Cells.Item(95, 11).Value = 101.0646450978375
$marketMaking.Cells.Item(95, 12).Value = 0
$marketMaking.Cells.Item(95, 13).Value = 0
$marketMaking.Cells.Item(95, 14).Value = 0.6
$marketMaking.Cells.Item(95, 15).Value = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(95, 17).Value = 0
